$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a trailing newline to the IP address values in column A (rows 2-6)
$ws.Range("A2").Value = "178.137.87.242`n"
$ws.Range("A3").Value = "46.148.22.18`n"
$ws.Range("A4").Value = "201.18.18.173`n"
$ws.Range("A5").Value = "46.148.18.162`n"
$ws.Range("A6").Value = "37.187.129.166`n"
